$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:D to B:E
$ws.Columns.Item(1).Insert()

# Copy header formatting (bold, border, centered) from B1 into new A1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set header text
$ws.Range("A1").Value = "ID"

# Populate ID values for rows 2-61
$ids = @("Hb 47","Hb 48","S 6","Hb 7","Hb 46","Hb 1","Hb 2","Hb 3","Hb 5","S 24","S 25","S 26","S 27","S 28","S 29","S 30","Hb 83","Hb 84","Hb 85","Hb 86","Hb 87","Hb 88","Hb 89","Hb 90","Hb 91","Hb 92","Hb 40","Hb 41","Hb 42","Hb 43","S 8","S 9","S 11","S 12","Hb 53","Hb 54","Hb 55","Hb 56","Hb 57","Hb 58","Hb 59","Hb 60","Hb 61","Hb 62","Hb 35","Hb 36","Hb 38","Hb 39","S 1","S 2","S 3","S 4","S 5","Hb 73","Hb 74","Hb 75","Hb 76","Hb 77","Hb 78","Hb 79")
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}

Write-Host "Inserted ID column with $($ids.Length) values"